# Refresh per-item market/profit figures in the Leve profit tracker.
# Each sheet (crafting class) lists leves with current market-board
# averages (H/I/J) and computed leve-turn-in price/profit (K/L/M/N).
# This batch re-syncs those columns with the latest pulled prices.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62: The Mustache Suits Him / Enchanted Mythrite Ink
$ws.Range("H62").Value = 31332.666
$ws.Range("I62").Value = 2999
$ws.Range("J62").Value = 45499.5
$ws.Range("K62").Value = 2999
$ws.Range("L62").Value = 45499.5
$ws.Range("M62").Value = -2375
$ws.Range("N62").Value = -46747.5

# Row 65: Forgery of Convenience (L) / Enchanted Mythrite Ink
$ws.Range("H65").Value = 31332.666
$ws.Range("I65").Value = 2999
$ws.Range("J65").Value = 45499.5
$ws.Range("K65").Value = 14995
$ws.Range("L65").Value = 227497.5
$ws.Range("M65").Value = -11875
$ws.Range("N65").Value = -233737.5

# Row 80: Cleansing the Wicked Humours / Hallowed Water
$ws.Range("H80").Value = 4716.6665
$ws.Range("I80").Value = 1900
$ws.Range("J80").Value = 5280
$ws.Range("K80").Value = 5700
$ws.Range("L80").Value = 15840
$ws.Range("M80").Value = -4702
$ws.Range("N80").Value = -17836

# Row 83: Washing Away the Sins (L) / Hallowed Water
$ws.Range("H83").Value = 4716.6665
$ws.Range("I83").Value = 1900
$ws.Range("J83").Value = 5280
$ws.Range("K83").Value = 17100
$ws.Range("L83").Value = 47520
$ws.Range("M83").Value = -12108
$ws.Range("N83").Value = -57504

# Row 88: The Grave of Hemlock Groves / Growth Formula Zeta
$ws.Range("H88").Value = 1011875.25
$ws.Range("J88").Value = 1347683.1
$ws.Range("L88").Value = 1347683.1
$ws.Range("N88").Value = -1348495.1

# Row 91: Dappling the Highlands (L) / Growth Formula Zeta
$ws.Range("H91").Value = 1011875.25
$ws.Range("J91").Value = 1347683.1
$ws.Range("L91").Value = 1347683.1
$ws.Range("N91").Value = -1350491.1

# Row 106: Making Your Mark / Enchanted Palladium Ink
$ws.Range("H106").Value = 10662.786
$ws.Range("I106").Value = 2409.75
$ws.Range("K106").Value = 2409.75
$ws.Range("M106").Value = -1778.75

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 3148.3
$ws.Range("I137").Value = 2122.516
$ws.Range("J137").Value = 4821.9473
$ws.Range("K137").Value = 6367.548000000001
$ws.Range("L137").Value = 14465.8419
$ws.Range("M137").Value = -3817.548000000001
$ws.Range("N137").Value = -19565.8419

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 3115.4138
$ws.Range("I138").Value = 1484.9166
$ws.Range("K138").Value = 4454.7498
$ws.Range("M138").Value = 685.2502000000004

# Row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 4686.5454
$ws.Range("I141").Value = 4394.1113
$ws.Range("J141").Value = 6002.5
$ws.Range("K141").Value = 13182.3339
$ws.Range("L141").Value = 18007.5
$ws.Range("M141").Value = -8002.333899999998
$ws.Range("N141").Value = -28367.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 10002432
$ws.Range("I32").Value = 10639863
$ws.Range("K32").Value = 10639863
$ws.Range("M32").Value = -10639576

# Row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("I45").Value = 55556332
$ws.Range("J45").Value = 2500
$ws.Range("K45").Value = 55556332
$ws.Range("L45").Value = 2500
$ws.Range("M45").Value = -55555955
$ws.Range("N45").Value = -3254

# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 13547420
$ws.Range("I61").Value = 18521632
$ws.Range("K61").Value = 18521632
$ws.Range("M61").Value = -18521420

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 6195208.5
$ws.Range("I74").Value = 10002857
$ws.Range("K74").Value = 10002857
$ws.Range("M74").Value = -10001983

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 6195208.5
$ws.Range("I77").Value = 10002857
$ws.Range("K77").Value = 50014285
$ws.Range("M77").Value = -50009917

# Row 102: Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws.Range("H102").Value = 8520.154
$ws.Range("I102").Value = 8522.817999999999
$ws.Range("K102").Value = 8522.817999999999
$ws.Range("M102").Value = -6900.817999999999

# Row 110: Scheduled Maintenance / Deepgold Ingot
$ws.Range("H110").Value = 1439.5
$ws.Range("I110").Value = 1439.5
$ws.Range("K110").Value = 1439.5
$ws.Range("M110").Value = 605.5

# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 5249.5
$ws.Range("I122").Value = 4000
$ws.Range("J122").Value = 5666
$ws.Range("K122").Value = 12000
$ws.Range("L122").Value = 16998
$ws.Range("M122").Value = -9550
$ws.Range("N122").Value = -21898

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 2706.1226
$ws.Range("I132").Value = 1415.641
$ws.Range("J132").Value = 7739
$ws.Range("K132").Value = 4246.923000000001
$ws.Range("L132").Value = 23217
$ws.Range("M132").Value = -1716.923000000001
$ws.Range("N132").Value = -28277

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 13547420
$ws.Range("I136").Value = 18521632
$ws.Range("K136").Value = 55564896
$ws.Range("M136").Value = -55562346

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 3112.6667
$ws.Range("I86").Value = 2419.1667
$ws.Range("K86").Value = 2419.1667
$ws.Range("M86").Value = -1296.1667

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 3112.6667
$ws.Range("I89").Value = 2419.1667
$ws.Range("K89").Value = 12095.8335
$ws.Range("M89").Value = -6479.833500000001

# Row 99: Meddle in Metal / Oroshigane Ingot
$ws.Range("H99").Value = 5919.069
$ws.Range("I99").Value = 7508.353
$ws.Range("K99").Value = 7508.353
$ws.Range("M99").Value = -6010.353

# Row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 1627.7273
$ws.Range("I105").Value = 1295
$ws.Range("K105").Value = 1295
$ws.Range("M105").Value = 452

# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 2096.4666
$ws.Range("I107").Value = 1678.1666
$ws.Range("J107").Value = 3769.6667
$ws.Range("K107").Value = 1678.1666
$ws.Range("L107").Value = 3769.6667
$ws.Range("M107").Value = 241.8334
$ws.Range("N107").Value = -7609.6667

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 401622.12
$ws.Range("I134").Value = 1416.1428
$ws.Range("J134").Value = 2502703.5
$ws.Range("K134").Value = 4248.428400000001
$ws.Range("L134").Value = 7508110.5
$ws.Range("M134").Value = -1713.428400000001
$ws.Range("N134").Value = -7513180.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 640022.4399999999
$ws.Range("I31").Value = 17678.334
$ws.Range("K31").Value = 17678.334
$ws.Range("M31").Value = -17383.334

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 640022.4399999999
$ws.Range("I34").Value = 17678.334
$ws.Range("K34").Value = 17678.334
$ws.Range("M34").Value = -17476.334

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 4067
$ws.Range("I132").Value = 4104.3125
$ws.Range("K132").Value = 12312.9375
$ws.Range("M132").Value = -9782.9375

$ws = $wb.Worksheets.Item("CUL")
# Row 139: Najoothie / Wild Banana Blend
$ws.Range("H139").Value = 2749.2593
$ws.Range("J139").Value = 2894.7368
$ws.Range("L139").Value = 8684.2104
$ws.Range("N139").Value = -18964.2104

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit / Mythrite Ingot
$ws.Range("H70").Value = 5533.4443
$ws.Range("I70").Value = 5257.2856
$ws.Range("K70").Value = 5257.2856
$ws.Range("M70").Value = -4987.2856

# Row 73: Hulls of Broken Dreams (L) / Mythrite Ingot
$ws.Range("H73").Value = 5533.4443
$ws.Range("I73").Value = 5257.2856
$ws.Range("K73").Value = 5257.2856
$ws.Range("M73").Value = -4321.2856

# Row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Range("H80").Value = 1849.5
$ws.Range("I80").Value = 1529.5714
$ws.Range("J80").Value = 2297.4
$ws.Range("K80").Value = 1529.5714
$ws.Range("L80").Value = 2297.4
$ws.Range("M80").Value = -531.5714
$ws.Range("N80").Value = -4293.4

# Row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("H83").Value = 1849.5
$ws.Range("I83").Value = 1529.5714
$ws.Range("J83").Value = 2297.4
$ws.Range("K83").Value = 7647.857
$ws.Range("L83").Value = 11487
$ws.Range("M83").Value = -2655.857
$ws.Range("N83").Value = -21471

# Row 95: Chain of Command / Koppranickel Temple Chain
$ws.Range("H95").Value = 55762.832
$ws.Range("J95").Value = 55762.832
$ws.Range("L95").Value = 55762.832
$ws.Range("N95").Value = -61254.832

# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 153.5
$ws.Range("I122").Value = 153.5
$ws.Range("K122").Value = 460.5
$ws.Range("M122").Value = 1989.5

# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 34488976
$ws.Range("I132").Value = 40002532
$ws.Range("K132").Value = 120007596
$ws.Range("M132").Value = -120005066

$ws = $wb.Worksheets.Item("LTW")
# Row 62: Pummeling Abroad / Archaeoskin Breeches of Maiming
$ws.Range("H62").Value = 48249
$ws.Range("J62").Value = 48249
$ws.Range("L62").Value = 48249
$ws.Range("N62").Value = -49497

# Row 65: The Style of the Time (L) / Archaeoskin Breeches of Maiming
$ws.Range("H65").Value = 48249
$ws.Range("J65").Value = 48249
$ws.Range("L65").Value = 144747
$ws.Range("N65").Value = -150987

# Row 101: A Stitch in Time / Marid Leather Gloves of Healing
$ws.Range("H101").Value = 100000
$ws.Range("J101").Value = 100000
$ws.Range("L101").Value = 100000
$ws.Range("N101").Value = -106490

# Row 133: The Perfect Accessory / Loboskin Amulet of Fending
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 43: Walk Softly and Carry a Big Halberd / Velveteen Dress Shoes
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()

# Row 103: To the Tops / Serge Gambison of Healing
$ws.Range("H103").Value = 43333.668
$ws.Range("J103").Value = 43333.668
$ws.Range("L103").Value = 43333.668
$ws.Range("N103").Value = -45677.668

# Row 107: Flax Wax / Bright Linen Yarn
$ws.Range("H107").Value = 18519612
$ws.Range("I107").Value = 33334832
$ws.Range("K107").Value = 100004496
$ws.Range("M107").Value = -100002576

# Row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 1793.6842
$ws.Range("I122").Value = 1912.6666
$ws.Range("K122").Value = 5737.9998
$ws.Range("M122").Value = -3287.9998

# Row 126: A Polished Purchase / Snow Linen
$ws.Range("H126").Value = 1341.0588
$ws.Range("I126").Value = 1341.0588
$ws.Range("K126").Value = 4023.1764
$ws.Range("M126").Value = -1553.1764

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 2190.697
$ws.Range("I132").Value = 1793.1
$ws.Range("K132").Value = 5379.299999999999
$ws.Range("M132").Value = -2849.299999999999

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 2126.85
$ws.Range("I136").Value = 1477.375
$ws.Range("K136").Value = 4432.125
$ws.Range("M136").Value = -1882.125
